# Weekly cryptos-list refresh (GitHub Actions bot).
# Updates the "Price" (D) and "Volume(1h)" (E) columns for each coin row
# with freshly scraped values. Columns D/E hold plain text (not numbers) in
# this sheet, so for values that look numeric we force the Text number
# format first -- otherwise Excel would silently reinterpret e.g. "0.9998"
# or "1.0000" as a number and drop the trailing zero / formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.339.30"
$ws.Range("E2").Value = "  +3.26%  "
$ws.Range("D3").Value = "2.006.23"
$ws.Range("E3").Value = "  +7.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7890"
$ws.Range("E5").Value = "  +67.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "258.72"
$ws.Range("E6").Value = "  +6.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3600"
$ws.Range("E8").Value = "  +25.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.68"
$ws.Range("E9").Value = "  +31.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07097"
$ws.Range("E10").Value = "  +9.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8598"
$ws.Range("E11").Value = "  +18.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08217"
$ws.Range("E12").Value = "  +5.39%  "
$ws.Range("D13").Value = "2.006.17"
$ws.Range("E13").Value = "  +7.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "101.46"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.618"
$ws.Range("E15").Value = "  +8.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "275.78"
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.94"
$ws.Range("E17").Value = "  +14.09%  "
$ws.Range("D18").Value = "31.348.15"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.916"
$ws.Range("E19").Value = "  +12.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000008000"
$ws.Range("E20").Value = "  +6.95%  "
$ws.Range("D21").Value = "2.269.51"
$ws.Range("E21").Value = "  +7.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.205"
$ws.Range("E24").Value = "  +14.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.16"
$ws.Range("E25").Value = "  +12.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1519"
$ws.Range("E26").Value = "  +57.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.06"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.11"
$ws.Range("E28").Value = "  +5.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.409"
$ws.Range("E29").Value = "  +27.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.623"
$ws.Range("E30").Value = "  +9.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.637"
$ws.Range("E31").Value = "  +9.46%  "
$ws.Range("E32").Value = "  +3.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.420"
$ws.Range("E33").Value = "  +6.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05256"
$ws.Range("E34").Value = "  +9.40%  "
$ws.Range("E35").Value = "  +9.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7726"
$ws.Range("E36").Value = "  +12.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.808"
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02018"
$ws.Range("E38").Value = "  +6.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.942"
$ws.Range("E39").Value = "  +3.79%  "
$ws.Range("E40").Value = "  +7.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "81.22"
$ws.Range("E41").Value = "  +7.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4772"
$ws.Range("E42").Value = "  +13.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.161"
$ws.Range("E43").Value = "  +10.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.60"
$ws.Range("E44").Value = "  +6.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8590"
$ws.Range("E45").Value = "  +4.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.802"
$ws.Range("E47").Value = "  +11.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.997"
$ws.Range("E48").Value = "  +3.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4380"
$ws.Range("E49").Value = "  +12.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.08"
$ws.Range("E50").Value = "  +5.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1190"
$ws.Range("E51").Value = "  +13.80%  "
